$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Duplicate "Nädal 9" to create the new "Nädal 10" sheet, placed
#    immediately after it (this also mirrors all of sheet 9's
#    formatting, merged cells, column widths, etc.).
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item("Nädal 9")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws10 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws10.Name = "Nädal 10"

# ------------------------------------------------------------------
# 2. Update the week-range header cell for the new week.
# ------------------------------------------------------------------
$ws10.Range("G4").Value = "31.03.2020 - 06.04.2020"

# ------------------------------------------------------------------
# 3. Replace the first (carried-over) log entry with the new week's
#    first entry: "48. osa - Quantity Facade testid ja klassid".
# ------------------------------------------------------------------
$ws10.Range("B7").Value = 43911
$ws10.Range("C7").Value = 0.51041666666666663
$ws10.Range("D7").Value = 0.56041666666666667
$ws10.Range("F7").Value = 72
$ws10.Range("G7").Value = "Kodutöö 9"
$ws10.Range("H7").Value = "p. 48 tehtud"

# ------------------------------------------------------------------
# 4. Clear out the remaining rows that were carried over from the
#    copy of "Nädal 9" so the new week starts blank beyond row 7
#    (rows 16-19 were already blank on the source sheet).
# ------------------------------------------------------------------
$ws10.Range("B8:J15").ClearContents()

# ------------------------------------------------------------------
# 5. Selection / active-tab bookkeeping to match the saved state.
#    Order matters: selecting a range on a sheet activates that
#    sheet, so we touch sheets 8 and 9 first and finish on sheet 10
#    so it ends up as the active tab.
# ------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("Nädal 8")
[void]$ws8.Range("G20").Select()

$ws9 = $wb.Worksheets.Item("Nädal 9")
[void]$ws9.Range("H19").Select()

[void]$ws10.Range("H17").Select()
